$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.957.05'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.796.26'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.37'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5300'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.36%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3975'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07479'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.30'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.086'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.211'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.75%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.522'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.34'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.794.40'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.54'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.73%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06604'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.02%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.17'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.972'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.996.05'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.06'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.42%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.091'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.44%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.71'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.18'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.003.08'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.16%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.316'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.26'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1096'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.090'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.676'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.54%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.500'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.50%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07036'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2214'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.92%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.166'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.94%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.37%  '

$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.27'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.00%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.382'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.192'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6123'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.83%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.406'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.49%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.40'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.676'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5721'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.47'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.94%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.182'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.922'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06802'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.46%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.20%  '
